$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05601
$ws.Range("N2").Value = 0.16803
$ws.Range("O2").Value = 0.02710547761971223
$ws.Range("P2").Value = 0.02710547761971223
$ws.Range("Q2").Value = 0.02005047847
$ws.Range("R2").Value = 0.18045430623
$ws.Range("S2").Value = 0.02710547761971223
$ws.Range("T2").Value = 0.02710547761971223

# Row 3 updates (precision-only changes)
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.902427539668559
$ws.Range("P3").Value = 0.9024275396685592
$ws.Range("Q3").Value = 0.6675441845636667
$ws.Range("R3").Value = 6.007897661073001
$ws.Range("S3").Value = 0.902427539668559
$ws.Range("T3").Value = 0.9024275396685592

# Row 4 updates
$ws.Range("M4").Value = 0.145611
$ws.Range("N4").Value = 0.436833
$ws.Range("O4").Value = 0.07046698271172858
$ws.Range("P4").Value = 0.07046698271172858
$ws.Range("Q4").Value = 0.05212587431700001
$ws.Range("R4").Value = 0.469132868853
$ws.Range("S4").Value = 0.07046698271172858
$ws.Range("T4").Value = 0.07046698271172858
